# Add isExpense field to the categories entity.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header for new column H ---
$ws.Cells.Item(1, 8).Value = "isExpense"

# --- Fill H2:H34 (all existing "expense" categories) with 1 ---
for ($r = 2; $r -le 34; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}

# --- Rows 35 & 36 ("דוח שנתי" section) also get isExpense = 1 ---
$ws.Cells.Item(35, 8).Value = 1
$ws.Cells.Item(36, 8).Value = 1

# --- New rows 37-41: income-type categories (isExpense = 0) ---
$newRows = @(
    @{ A = "הכנסה מעסק"; B = "הכנסה מעסק" },
    @{ A = "משכורת";     B = "משכורת" },
    @{ A = "ביטוח לאומי"; B = "קצבת ילדים" },
    @{ A = "ביטוח לאומי"; B = "מילואים" },
    @{ A = "ביטוח לאומי"; B = "דמי לידה" }
)

$r = 37
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
    $ws.Cells.Item($r, 5).Value = 0
    $ws.Cells.Item($r, 6).Value = 0
    $ws.Cells.Item($r, 7).Value = 0
    $ws.Cells.Item($r, 8).Value = 0
    $r++
}

# --- Column widths: widen G, add width for new column H ---
$ws.Range("G1").EntireColumn.ColumnWidth = 27.796875
$ws.Range("H1").EntireColumn.ColumnWidth = 21.3984375

# --- View: scroll position & selection ---
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D43").Select()
